$wb = $excel.ActiveWorkbook

# --- network_weights: tiny precision correction on D5 (recalculated value) ---
$wsWeights = $wb.Worksheets.Item("network_weights")
$wsWeights.Range("D5").Value2 = -0.97501548477922195
$wsWeights.Activate()
$wsWeights.Range("K34").Select()

# --- optimization_parameters: remove the stray leftover "Sheet" row (old row 16) ---
$wsParams = $wb.Worksheets.Item("optimization_parameters")
$wsParams.Rows.Item(16).Delete()
$wsParams.Activate()
$wsParams.Rows.Item(16).Select()

# --- optimization_diagnostics becomes the active sheet/tab ---
$wsDiag = $wb.Worksheets.Item("optimization_diagnostics")
$wsDiag.Activate()
